$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 5111.1113
$ws.Cells.Item(12, 9).Value = 9000
$ws.Cells.Item(12, 10).Value = 2000
$ws.Cells.Item(12, 11).Value = 9000
$ws.Cells.Item(12, 12).Value = 2000
$ws.Cells.Item(12, 13).Value = -8830
$ws.Cells.Item(12, 14).Value = -2340
$ws.Cells.Item(17, 8).Value = 557306.4399999999
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 557306.4399999999
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 1671919.32
$ws.Cells.Item(17, 13).ClearContents()
$ws.Cells.Item(17, 14).Value = -1672255.32
$ws.Cells.Item(38, 8).Value = 552.6
$ws.Cells.Item(38, 9).Value = 582.2727
$ws.Cells.Item(38, 10).Value = 471
$ws.Cells.Item(38, 11).Value = 1746.8181
$ws.Cells.Item(38, 12).Value = 1413
$ws.Cells.Item(38, 13).Value = -1374.8181
$ws.Cells.Item(38, 14).Value = -2157
$ws.Cells.Item(41, 8).Value = 745.625
$ws.Cells.Item(41, 10).Value = 360.16666
$ws.Cells.Item(41, 12).Value = 360.16666
$ws.Cells.Item(41, 14).Value = -1240.16666
$ws.Cells.Item(47, 8).Value = 8922.333000000001
$ws.Cells.Item(47, 9).Value = 8922.333000000001
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 11).Value = 8922.333000000001
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 13).Value = -7950.333000000001
$ws.Cells.Item(47, 14).ClearContents()
$ws.Cells.Item(55, 8).Value = 553
$ws.Cells.Item(55, 10).Value = 626
$ws.Cells.Item(55, 12).Value = 626
$ws.Cells.Item(55, 14).Value = -1054
$ws.Cells.Item(76, 8).Value = 3562.0715
$ws.Cells.Item(76, 9).Value = 3570.0908
$ws.Cells.Item(76, 10).Value = 3532.6667
$ws.Cells.Item(76, 11).Value = 3570.0908
$ws.Cells.Item(76, 12).Value = 3532.6667
$ws.Cells.Item(76, 13).Value = -3255.0908
$ws.Cells.Item(76, 14).Value = -4162.6667
$ws.Cells.Item(79, 8).Value = 3562.0715
$ws.Cells.Item(79, 9).Value = 3570.0908
$ws.Cells.Item(79, 10).Value = 3532.6667
$ws.Cells.Item(79, 11).Value = 3570.0908
$ws.Cells.Item(79, 12).Value = 3532.6667
$ws.Cells.Item(79, 13).Value = -2478.0908
$ws.Cells.Item(79, 14).Value = -5716.6667
$ws.Cells.Item(86, 8).Value = 931.6667
$ws.Cells.Item(86, 9).Value = 931.6667
$ws.Cells.Item(86, 11).Value = 931.6667
$ws.Cells.Item(86, 13).Value = 191.3333
$ws.Cells.Item(89, 8).Value = 931.6667
$ws.Cells.Item(89, 9).Value = 931.6667
$ws.Cells.Item(89, 11).Value = 4658.3335
$ws.Cells.Item(89, 13).Value = 957.6665000000003
$ws.Cells.Item(92, 8).Value = 612.6316
$ws.Cells.Item(92, 9).Value = 543.7857
$ws.Cells.Item(92, 11).Value = 543.7857
$ws.Cells.Item(92, 13).Value = 704.2143
$ws.Cells.Item(137, 8).Value = 1062.72
$ws.Cells.Item(137, 9).Value = 992.4211
$ws.Cells.Item(137, 10).Value = 1285.3334
$ws.Cells.Item(137, 11).Value = 2977.2633
$ws.Cells.Item(137, 12).Value = 3856.0002
$ws.Cells.Item(137, 13).Value = -427.2633000000001
$ws.Cells.Item(137, 14).Value = -8956.0002
$ws.Cells.Item(139, 8).Value = 149999.83
$ws.Cells.Item(139, 10).Value = 149999.83
$ws.Cells.Item(139, 12).Value = 149999.83
$ws.Cells.Item(139, 14).Value = -160279.83

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 6995.577
$ws.Cells.Item(45, 9).Value = 8165.6665
$ws.Cells.Item(45, 10).Value = 4362.875
$ws.Cells.Item(45, 11).Value = 8165.6665
$ws.Cells.Item(45, 12).Value = 4362.875
$ws.Cells.Item(45, 13).Value = -7788.6665
$ws.Cells.Item(45, 14).Value = -5116.875
$ws.Cells.Item(57, 8).Value = 5168.5713
$ws.Cells.Item(57, 9).Value = 5168.5713
$ws.Cells.Item(57, 11).Value = 5168.5713
$ws.Cells.Item(57, 13).Value = -4684.5713
$ws.Cells.Item(61, 8).Value = 1559.591
$ws.Cells.Item(61, 9).Value = 1272.8667
$ws.Cells.Item(61, 10).Value = 2174
$ws.Cells.Item(61, 11).Value = 1272.8667
$ws.Cells.Item(61, 12).Value = 2174
$ws.Cells.Item(61, 13).Value = -1060.8667
$ws.Cells.Item(61, 14).Value = -2598
$ws.Cells.Item(74, 8).Value = 1506.8077
$ws.Cells.Item(74, 9).Value = 1553.25
$ws.Cells.Item(74, 10).Value = 949.5
$ws.Cells.Item(74, 11).Value = 1553.25
$ws.Cells.Item(74, 12).Value = 949.5
$ws.Cells.Item(74, 13).Value = -679.25
$ws.Cells.Item(74, 14).Value = -2697.5
$ws.Cells.Item(77, 8).Value = 1506.8077
$ws.Cells.Item(77, 9).Value = 1553.25
$ws.Cells.Item(77, 10).Value = 949.5
$ws.Cells.Item(77, 11).Value = 7766.25
$ws.Cells.Item(77, 12).Value = 4747.5
$ws.Cells.Item(77, 13).Value = -3398.25
$ws.Cells.Item(77, 14).Value = -13483.5
$ws.Cells.Item(136, 8).Value = 1559.591
$ws.Cells.Item(136, 9).Value = 1272.8667
$ws.Cells.Item(136, 10).Value = 2174
$ws.Cells.Item(136, 11).Value = 3818.6001
$ws.Cells.Item(136, 12).Value = 6522
$ws.Cells.Item(136, 13).Value = -1268.6001
$ws.Cells.Item(136, 14).Value = -11622

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(113, 8).Value = 5126.6665
$ws.Cells.Item(113, 9).Value = 5126.6665
$ws.Cells.Item(113, 11).Value = 5126.6665
$ws.Cells.Item(113, 13).Value = -2956.6665
$ws.Cells.Item(128, 8).Value = 3461.8333
$ws.Cells.Item(128, 9).Value = 3461.8333
$ws.Cells.Item(128, 11).Value = 10385.4999
$ws.Cells.Item(128, 13).Value = -7895.499899999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 8930.842000000001
$ws.Cells.Item(31, 9).Value = 2789.625
$ws.Cells.Item(31, 10).Value = 23380.766
$ws.Cells.Item(31, 11).Value = 2789.625
$ws.Cells.Item(31, 12).Value = 23380.766
$ws.Cells.Item(31, 13).Value = -2494.625
$ws.Cells.Item(31, 14).Value = -23970.766
$ws.Cells.Item(34, 8).Value = 8930.842000000001
$ws.Cells.Item(34, 9).Value = 2789.625
$ws.Cells.Item(34, 10).Value = 23380.766
$ws.Cells.Item(34, 11).Value = 2789.625
$ws.Cells.Item(34, 12).Value = 23380.766
$ws.Cells.Item(34, 13).Value = -2587.625
$ws.Cells.Item(34, 14).Value = -23784.766
$ws.Cells.Item(58, 8).Value = 2049.75
$ws.Cells.Item(58, 9).Value = 1149.5
$ws.Cells.Item(58, 10).Value = 2950
$ws.Cells.Item(58, 11).Value = 1149.5
$ws.Cells.Item(58, 12).Value = 2950
$ws.Cells.Item(58, 13).Value = -946.5
$ws.Cells.Item(58, 14).Value = -3356
$ws.Cells.Item(68, 8).Value = 24971.428
$ws.Cells.Item(68, 10).Value = 24971.428
$ws.Cells.Item(68, 12).Value = 24971.428
$ws.Cells.Item(68, 14).Value = -26469.428
$ws.Cells.Item(71, 8).Value = 24971.428
$ws.Cells.Item(71, 10).Value = 24971.428
$ws.Cells.Item(71, 12).Value = 74914.284
$ws.Cells.Item(71, 14).Value = -82402.284
$ws.Cells.Item(74, 8).Value = 57250
$ws.Cells.Item(74, 10).Value = 57250
$ws.Cells.Item(74, 12).Value = 57250
$ws.Cells.Item(74, 14).Value = -58998
$ws.Cells.Item(76, 8).Value = 5073.3335
$ws.Cells.Item(76, 9).Value = 5073.3335
$ws.Cells.Item(76, 11).Value = 5073.3335
$ws.Cells.Item(76, 13).Value = -4758.3335
$ws.Cells.Item(77, 8).Value = 57250
$ws.Cells.Item(77, 10).Value = 57250
$ws.Cells.Item(77, 12).Value = 171750
$ws.Cells.Item(77, 14).Value = -180486
$ws.Cells.Item(79, 8).Value = 5073.3335
$ws.Cells.Item(79, 9).Value = 5073.3335
$ws.Cells.Item(79, 11).Value = 5073.3335
$ws.Cells.Item(79, 13).Value = -3981.3335
$ws.Cells.Item(136, 8).Value = 2049.75
$ws.Cells.Item(136, 9).Value = 1149.5
$ws.Cells.Item(136, 10).Value = 2950
$ws.Cells.Item(136, 11).Value = 3448.5
$ws.Cells.Item(136, 12).Value = 8850
$ws.Cells.Item(136, 13).Value = -898.5
$ws.Cells.Item(136, 14).Value = -13950

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 420.26666
$ws.Cells.Item(2, 9).Value = 1415.75
$ws.Cells.Item(2, 10).Value = 58.272728
$ws.Cells.Item(2, 11).Value = 8494.5
$ws.Cells.Item(2, 12).Value = 349.636368
$ws.Cells.Item(2, 13).Value = -8381.5
$ws.Cells.Item(2, 14).Value = -575.6363679999999
$ws.Cells.Item(4, 8).Value = 6176833.5
$ws.Cells.Item(4, 10).Value = 4882533.5
$ws.Cells.Item(4, 12).Value = 14647600.5
$ws.Cells.Item(4, 14).Value = -14647824.5
$ws.Cells.Item(7, 8).Value = 65.666664
$ws.Cells.Item(7, 10).Value = 97.5
$ws.Cells.Item(7, 12).Value = 292.5
$ws.Cells.Item(7, 14).Value = -516.5
$ws.Cells.Item(11, 8).Value = 1364218.5
$ws.Cells.Item(11, 9).Value = 2143004.2
$ws.Cells.Item(11, 10).Value = 1343.5
$ws.Cells.Item(11, 11).Value = 6429012.600000001
$ws.Cells.Item(11, 12).Value = 4030.5
$ws.Cells.Item(11, 13).Value = -6428872.600000001
$ws.Cells.Item(11, 14).Value = -4310.5
$ws.Cells.Item(62, 8).Value = 11499.5
$ws.Cells.Item(62, 10).Value = 11499.5
$ws.Cells.Item(62, 12).Value = 34498.5
$ws.Cells.Item(62, 14).Value = -35870.5
$ws.Cells.Item(65, 8).Value = 11499.5
$ws.Cells.Item(65, 10).Value = 11499.5
$ws.Cells.Item(65, 12).Value = 103495.5
$ws.Cells.Item(65, 14).Value = -110359.5
$ws.Cells.Item(68, 8).Value = 1355.1428
$ws.Cells.Item(68, 9).Value = 1374.5
$ws.Cells.Item(68, 11).Value = 4123.5
$ws.Cells.Item(68, 13).Value = -3312.5
$ws.Cells.Item(71, 8).Value = 1355.1428
$ws.Cells.Item(71, 9).Value = 1374.5
$ws.Cells.Item(71, 11).Value = 12370.5
$ws.Cells.Item(71, 13).Value = -8314.5
$ws.Cells.Item(86, 8).Value = 588.8946999999999
$ws.Cells.Item(86, 9).Value = 485.1
$ws.Cells.Item(86, 11).Value = 1455.3
$ws.Cells.Item(86, 13).Value = -269.3000000000002
$ws.Cells.Item(89, 8).Value = 588.8946999999999
$ws.Cells.Item(89, 9).Value = 485.1
$ws.Cells.Item(89, 11).Value = 4365.900000000001
$ws.Cells.Item(89, 13).Value = 1562.099999999999
$ws.Cells.Item(103, 8).Value = 633.0909
$ws.Cells.Item(103, 9).Value = 510.83334
$ws.Cells.Item(103, 10).Value = 779.8
$ws.Cells.Item(103, 11).Value = 1532.50002
$ws.Cells.Item(103, 12).Value = 2339.4
$ws.Cells.Item(103, 13).Value = -653.5000199999999
$ws.Cells.Item(103, 14).Value = -4097.4
$ws.Cells.Item(131, 8).Value = 2257.8
$ws.Cells.Item(131, 10).Value = 1753
$ws.Cells.Item(131, 12).Value = 5259
$ws.Cells.Item(131, 14).Value = -15339
$ws.Cells.Item(132, 8).Value = 1687.0385
$ws.Cells.Item(132, 9).Value = 879.5454999999999
$ws.Cells.Item(132, 11).Value = 7915.9095
$ws.Cells.Item(132, 13).Value = -5385.9095

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(63, 8).Value = 20335.5
$ws.Cells.Item(63, 10).Value = 20335.5
$ws.Cells.Item(63, 12).Value = 20335.5
$ws.Cells.Item(63, 14).Value = -21707.5
$ws.Cells.Item(66, 8).Value = 20335.5
$ws.Cells.Item(66, 10).Value = 20335.5
$ws.Cells.Item(66, 12).Value = 61006.5
$ws.Cells.Item(66, 14).Value = -67870.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 5980
$ws.Cells.Item(136, 9).Value = 4995
$ws.Cells.Item(136, 10).Value = 6177
$ws.Cells.Item(136, 11).Value = 14985
$ws.Cells.Item(136, 12).Value = 18531
$ws.Cells.Item(136, 13).Value = -12435
$ws.Cells.Item(136, 14).Value = -23631

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1850.9688
$ws.Cells.Item(132, 9).Value = 1336.05
$ws.Cells.Item(132, 10).Value = 2709.1667
$ws.Cells.Item(132, 11).Value = 4008.15
$ws.Cells.Item(132, 12).Value = 8127.500100000001
$ws.Cells.Item(132, 13).Value = -1478.15
$ws.Cells.Item(132, 14).Value = -13187.5001
